$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously held one placeholder row (A2:B2 = 47323424/47323424).
# Replace it with the parsed order-line rows (source: JSON -> xlsx import).
$rows = @(
    @(357469, "0000006274", "Крем-борщ вегетарианский", 1, 127),
    @(357470, "0000006273", "Окрошка на квасе", 1, 135),
    @(357471, "0000006255", "Оливье с курой", 1, 151),
    @(357473, "0000002843", "Рис с курицей в кисло-сладком соусе", 1, 151),
    @(357477, "0000006260", "Бефстроганов с жареным картофелем и грибами", 1, 183),
    @(357478, "0000002878", "SWEETBOX чиа-манго", 1, 95)
)

$firstRow = 2
$lastRow = $firstRow + $rows.Length - 1

$r = $firstRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    # Column B holds zero-padded numeric-looking codes ("0000006274") that
    # must stay text. A direct .Value assignment gets auto-coerced to a
    # number (dropping the leading zeros), so write it as a `="..."` text
    # formula first and bake it down to a literal value below.
    $ws.Cells.Item($r, 2).Formula = '="' + $row[1] + '"'
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# Collapse the column-B text formulas down to plain literal (shared-string)
# values, same end state as typing the text directly into a text-formatted
# cell, without leaving a number-format override behind.
$codeRange = $ws.Range("B" + $firstRow + ":B" + $lastRow)
$codeRange.Copy()
$codeRange.PasteSpecial(-4163)
